# Generate Report for Handback
# Marks the zh-cn / de-de localization rows as handed back: updates the
# status text, populates the "Latest Target File" / "Latest Handback File"
# / "Latest Handback DateTime" columns, and adds a hyperlink on the new
# target-file cell (mirroring the existing source-file hyperlink).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceFileDisplay = "c6159970-c0da-4760-9c8e-2a9162e7e16e.md"
$sourceFileUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e8f5bc3f05499ddc35523817eb3838ef0eb39260/e2e/c6159970-c0da-4760-9c8e-2a9162e7e16e.md"

# ---------------------------------------------------------------------
# Overview sheet: status shown for both locales flips to "handed back"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

# Widen the zh-cn / de-de status columns to fit the longer handed-back text
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# Helper: stamp one locale sheet (zh-cn / de-de) as handed back
# ---------------------------------------------------------------------
function Set-HandbackRow($ws, $handbackDateTime) {
    $ws.Range("C2").Value = $statusText

    $targetFile = $ws.Range("G2").Text

    $ws.Hyperlinks.Add($ws.Range("I2"), $sourceFileUrl, $null, $null, $sourceFileDisplay) | Out-Null
    $ws.Range("I2").Font.Underline = 2
    $ws.Range("I2").Font.Color = 15570276

    $ws.Range("J2").Value = $targetFile
    $ws.Range("K2").Value = $handbackDateTime

    # Widen Status / Latest Target File / Latest Handback File columns to
    # fit the new, longer cell contents.
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackRow $wsZhCn "2016-08-23 12:59:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackRow $wsDeDe "2016-08-23 12:59:53"
